$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had three region rows (1-3) followed by a single
# "   " total row (row 4) that summed B1:B3. The edit turns that single
# trailing total row into a per-region "render label" row placed right
# after each region's own row (pushing the rolled-up label to the back
# of each region instead of only at the very end), each one referencing
# the value directly above it instead of summing everything.

# Make room: insert a blank row right after region A (row 1) and right
# after region B (which will be row 3 once the first insert shifts it
# down from row 2).
$ws.Rows("2:2").Insert()
$ws.Rows("4:4").Insert()

# New row 2 - label row for region A, referencing its value in B1.
$ws.Range("A2").Value = "   "
$ws.Range("B2").Formula = "=B1"

# New row 4 - label row for region B, referencing its value in B3.
$ws.Range("A4").Value = "   "
$ws.Range("B4").Formula = "=B3"

# The old trailing total row (originally row 4, now pushed down to row 6)
# becomes the label row for region C, referencing its value in B5 instead
# of summing the whole column.
$ws.Range("A6").Value = "   "
$ws.Range("B6").Formula = "=B5"

# Highlight the three original region value cells with the themed fill
# (Accent 6, Lighter 40%) so they stand out from their label rows below.
$ws.Range("B1").Interior.ThemeColor = 10
$ws.Range("B1").Interior.TintAndShade = 0.39997558519241921
$ws.Range("B1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the active selection on the first new label row, as in the source.
$ws.Range("A4").Select()
